$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-stamp A216:A280 with the existing date/time style (style index 1, already used
# throughout column A) so every row -- including the brand-new ones beyond the old
# A1:C243 extent -- carries "s=\"1\"" on its A cell, matching the pattern used by the
# sheet's other (still-blank) rows further down before this edit.
$ws.Range("A216:A280").NumberFormat = "dd/mm/yyyy\ h:mm:ss"

$data = New-Object 'object[,]' 64,3
$data[0,0] = 45071.40763883102
$data[0,1] = 1
$data[0,2] = 2
$data[1,0] = 45071.407986053244
$data[1,1] = 1
$data[1,2] = 3
$data[2,0] = 45071.40833327546
$data[2,1] = 1
$data[2,2] = 2
$data[3,0] = 45071.408680497683
$data[3,1] = 1
$data[3,2] = 3
$data[4,0] = 45071.409027719907
$data[4,1] = 1
$data[4,2] = 2
$data[5,0] = 45071.40937494213
$data[5,1] = 1
$data[5,2] = 3
$data[6,0] = 45071.409722164353
$data[6,1] = 1
$data[6,2] = 2
$data[7,0] = 45071.410069386577
$data[7,1] = 1
$data[7,2] = 3
$data[8,0] = 45071.410416608793
$data[8,1] = 1
$data[8,2] = 2
$data[9,0] = 45071.410763831016
$data[9,1] = 1
$data[9,2] = 3
$data[10,0] = 45071.411111053239
$data[10,1] = 1
$data[10,2] = 1
$data[11,0] = 45071.411458275463
$data[11,1] = 1
$data[11,2] = 1
$data[12,0] = 45071.411805497686
$data[12,1] = 1
$data[12,2] = 1
$data[13,0] = 45071.412152719909
$data[13,1] = 1
$data[13,2] = 1
$data[14,0] = 45071.41249988426
$data[14,1] = 1
$data[14,2] = 3
$data[15,0] = 45071.412847106483
$data[15,1] = 1
$data[15,2] = 2
$data[16,0] = 45071.413194328707
$data[16,1] = 1
$data[16,2] = 3
$data[17,0] = 45071.413541550923
$data[17,1] = 1
$data[17,2] = 2
$data[18,0] = 45071.413888773146
$data[18,1] = 1
$data[18,2] = 3
$data[19,0] = 45071.414235995369
$data[19,1] = 1
$data[19,2] = 2
$data[20,0] = 45071.414583217593
$data[20,1] = 1
$data[20,2] = 3
$data[21,0] = 45071.414930439816
$data[21,1] = 1
$data[21,2] = 1
$data[22,0] = 45071.415277662039
$data[22,1] = 1
$data[22,2] = 1
$data[23,0] = 45071.415624884263
$data[23,1] = 1
$data[23,2] = 1
$data[24,0] = 45071.415972106479
$data[24,1] = 1
$data[24,2] = 1
$data[25,0] = 45071.416319328702
$data[25,1] = 1
$data[25,2] = 3
$data[26,0] = 45071.416666550926
$data[26,1] = 1
$data[26,2] = 2
$data[27,0] = 45071.417013773149
$data[27,1] = 1
$data[27,2] = 3
$data[28,0] = 45071.417360995372
$data[28,1] = 1
$data[28,2] = 1
$data[29,0] = 45071.417708217596
$data[29,1] = 1
$data[29,2] = 1
$data[30,0] = 45071.418055439812
$data[30,1] = 1
$data[30,2] = 1
$data[31,0] = 45071.418402662035
$data[31,1] = 1
$data[31,2] = 1
$data[32,0] = 45071.418749884258
$data[32,1] = 1
$data[32,2] = 1
$data[33,0] = 45071.419097106482
$data[33,1] = 1
$data[33,2] = 1
$data[34,0] = 45071.419444328705
$data[34,1] = 1
$data[34,2] = 1
$data[35,0] = 45071.419791550928
$data[35,1] = 1
$data[35,2] = 1
$data[36,0] = 45071.420138773145
$data[36,1] = 1
$data[36,2] = 3
$data[37,0] = 45071.420485995368
$data[37,1] = 1
$data[37,2] = 2
$data[38,0] = 45071.420833217591
$data[38,1] = 1
$data[38,2] = 2
$data[39,0] = 45071.421180439815
$data[39,1] = 1
$data[39,2] = 2
$data[40,0] = 45071.421527662038
$data[40,1] = 1
$data[40,2] = 2
$data[41,0] = 45071.421874826388
$data[41,1] = 1
$data[41,2] = 2
$data[42,0] = 45071.422222048612
$data[42,1] = 1
$data[42,2] = 3
$data[43,0] = 45071.422569270835
$data[43,1] = 1
$data[43,2] = 1
$data[44,0] = 45071.422916493058
$data[44,1] = 1
$data[44,2] = 1
$data[45,0] = 45071.423263715275
$data[45,1] = 1
$data[45,2] = 1
$data[46,0] = 45071.423610937498
$data[46,1] = 1
$data[46,2] = 1
$data[47,0] = 45071.423958159721
$data[47,1] = 1
$data[47,2] = 1
$data[48,0] = 45071.424305381945
$data[48,1] = 1
$data[48,2] = 1
$data[49,0] = 45071.424652604168
$data[49,1] = 1
$data[49,2] = 3
$data[50,0] = 45071.424999826391
$data[50,1] = 1
$data[50,2] = 2
$data[51,0] = 45071.425347048615
$data[51,1] = 1
$data[51,2] = 3
$data[52,0] = 45071.425694270831
$data[52,1] = 1
$data[52,2] = 2
$data[53,0] = 45071.426041493054
$data[53,1] = 1
$data[53,2] = 2
$data[54,0] = 45071.426388715277
$data[54,1] = 1
$data[54,2] = 2
$data[55,0] = 45071.426735937501
$data[55,1] = 1
$data[55,2] = 2
$data[56,0] = 45071.427083159724
$data[56,1] = 1
$data[56,2] = 2
$data[57,0] = 45071.427430381947
$data[57,1] = 1
$data[57,2] = 2
$data[58,0] = 45071.427777604164
$data[58,1] = 1
$data[58,2] = 2
$data[59,0] = 45071.428124826387
$data[59,1] = 1
$data[59,2] = 2
$data[60,0] = 45071.42847204861
$data[60,1] = 1
$data[60,2] = 2
$data[61,0] = 45071.428819270834
$data[61,1] = 1
$data[61,2] = 2
$data[62,0] = 45071.429166493057
$data[62,1] = 1
$data[62,2] = 2
$data[63,0] = 45071.42951371528
$data[63,1] = 1
$data[63,2] = 2

$ws.Range("A216:C279").Value = $data

# Row 280 stays an empty placeholder (style only, no value) -- same "blank but styled"
# pattern the sheet already used for its trailing rows.

# Restore the view state: selection on C257, scrolled so row 233 is at the top.
$ws.Range("C257").Select()
$excel.ActiveWindow.ScrollRow = 233
